$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 778.8333
$ws.Range("I6").Value = 334.4
$ws.Range("K6").Value = 1003.2
$ws.Range("M6").Value = -891.1999999999999
$ws.Range("H32").Value = 4047.0417
$ws.Range("I32").Value = 6595.8
$ws.Range("K32").Value = 6595.8
$ws.Range("M32").Value = -6269.8
$ws.Range("H40").Value = 4726.5625
$ws.Range("I40").Value = 3342.4
$ws.Range("J40").Value = 7033.5
$ws.Range("K40").Value = 3342.4
$ws.Range("L40").Value = 7033.5
$ws.Range("M40").Value = -3167.4
$ws.Range("N40").Value = -7383.5
$ws.Range("H64").Value = 8179.5
$ws.Range("J64").Value = 11114.286
$ws.Range("L64").Value = 11114.286
$ws.Range("N64").Value = -11610.286
$ws.Range("H67").Value = 8179.5
$ws.Range("J67").Value = 11114.286
$ws.Range("L67").Value = 11114.286
$ws.Range("N67").Value = -12830.286
$ws.Range("H98").Value = 6704.1333
$ws.Range("I98").Value = 6684.0713
$ws.Range("K98").Value = 6684.0713
$ws.Range("M98").Value = -5186.0713
$ws.Range("H99").Value = 242.66667
$ws.Range("I99").Value = 242.66667
$ws.Range("K99").Value = 728.00001
$ws.Range("M99").Value = 769.99999
$ws.Range("H107").Value = 933.4286
$ws.Range("I107").Value = 954.6875
$ws.Range("K107").Value = 954.6875
$ws.Range("M107").Value = 965.3125
$ws.Range("H122").Value = 6704.1333
$ws.Range("I122").Value = 6684.0713
$ws.Range("K122").Value = 20052.2139
$ws.Range("M122").Value = -17602.2139
$ws.Range("H131").Value = 114743.22
$ws.Range("I131").Value = 114743.22
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 344229.66
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -339189.66
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 6648.7896
$ws.Range("I132").Value = 7570.0625
$ws.Range("K132").Value = 22710.1875
$ws.Range("M132").Value = -20180.1875
$ws.Range("H137").Value = 1426.24
$ws.Range("I137").Value = 1237.6666
$ws.Range("K137").Value = 3712.9998
$ws.Range("M137").Value = -1162.9998
$ws.Range("H141").Value = 9384.23
$ws.Range("I141").Value = 7700
$ws.Range("J141").Value = 14998.333
$ws.Range("K141").Value = 23100
$ws.Range("L141").Value = 44994.999
$ws.Range("M141").Value = -17920
$ws.Range("N141").Value = -55354.999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1232
$ws.Range("H32").Value = 3273.0164
$ws.Range("I32").Value = 3160.9
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 3160.9
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -2873.9
$ws.Range("N32").Value = -10574
$ws.Range("H45").Value = 62068.887
$ws.Range("I45").Value = 98783.52
$ws.Range("K45").Value = 98783.52
$ws.Range("M45").Value = -98406.52
$ws.Range("H61").Value = 4392470.5
$ws.Range("I61").Value = 5055943
$ws.Range("K61").Value = 5055943
$ws.Range("M61").Value = -5055731
$ws.Range("H122").Value = 1399.6923
$ws.Range("I122").Value = 1433
$ws.Range("K122").Value = 4299
$ws.Range("M122").Value = -1849
$ws.Range("H136").Value = 4392470.5
$ws.Range("I136").Value = 5055943
$ws.Range("K136").Value = 15167829
$ws.Range("M136").Value = -15165279
$ws.Range("H137").Value = 69993
$ws.Range("J137").Value = 69993
$ws.Range("L137").Value = 69993
$ws.Range("N137").Value = -80193
$ws.Range("H139").Value = 165753.8
$ws.Range("J139").Value = 165753.8
$ws.Range("L139").Value = 165753.8
$ws.Range("N139").Value = -176033.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4392.722
$ws.Range("I107").Value = 3672.8333
$ws.Range("J107").Value = 5832.5
$ws.Range("K107").Value = 3672.8333
$ws.Range("L107").Value = 5832.5
$ws.Range("M107").Value = -1752.8333
$ws.Range("N107").Value = -9672.5
$ws.Range("H134").Value = 5762.5
$ws.Range("I134").Value = 5878.75
$ws.Range("K134").Value = 17636.25
$ws.Range("M134").Value = -15101.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 692.6
$ws.Range("I6").Value = 692.6
$ws.Range("K6").Value = 692.6
$ws.Range("M6").Value = -579.6
$ws.Range("H31").Value = 6072.615
$ws.Range("I31").Value = 5871.6665
$ws.Range("J31").Value = 6244.857
$ws.Range("K31").Value = 5871.6665
$ws.Range("L31").Value = 6244.857
$ws.Range("M31").Value = -5576.6665
$ws.Range("N31").Value = -6834.857
$ws.Range("H34").Value = 6072.615
$ws.Range("I34").Value = 5871.6665
$ws.Range("J34").Value = 6244.857
$ws.Range("K34").Value = 5871.6665
$ws.Range("L34").Value = 6244.857
$ws.Range("M34").Value = -5669.6665
$ws.Range("N34").Value = -6648.857
$ws.Range("H58").Value = 6900.923
$ws.Range("I58").Value = 2652
$ws.Range("K58").Value = 2652
$ws.Range("M58").Value = -2449
$ws.Range("H86").Value = 8149
$ws.Range("I86").Value = 7298.3335
$ws.Range("J86").Value = 8999.666999999999
$ws.Range("K86").Value = 7298.3335
$ws.Range("L86").Value = 8999.666999999999
$ws.Range("M86").Value = -6175.3335
$ws.Range("N86").Value = -11245.667
$ws.Range("H89").Value = 8149
$ws.Range("I89").Value = 7298.3335
$ws.Range("J89").Value = 8999.666999999999
$ws.Range("K89").Value = 36491.6675
$ws.Range("L89").Value = 44998.335
$ws.Range("M89").Value = -30875.6675
$ws.Range("N89").Value = -56230.335
$ws.Range("H132").Value = 1285.091
$ws.Range("I132").Value = 1215.1111
$ws.Range("K132").Value = 3645.3333
$ws.Range("M132").Value = -1115.3333
$ws.Range("H136").Value = 6900.923
$ws.Range("I136").Value = 2652
$ws.Range("K136").Value = 7956
$ws.Range("M136").Value = -5406
$ws.Range("H141").Value = 37888
$ws.Range("J141").Value = 37888
$ws.Range("L141").Value = 37888
$ws.Range("N141").Value = -48248

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 33719812
$ws.Range("I4").Value = 32173446
$ws.Range("K4").Value = 96520338
$ws.Range("M4").Value = -96520226
$ws.Range("H99").Value = 1750
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H139").Value = 2186.9119
$ws.Range("I139").Value = 1619.8928
$ws.Range("K139").Value = 4859.678400000001
$ws.Range("M139").Value = 280.3215999999993

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1473914.1
$ws.Range("I3").Value = 670499.7
$ws.Range("K3").Value = 670499.7
$ws.Range("M3").Value = -670383.7
$ws.Range("H46").Value = 35419.5
$ws.Range("J46").Value = 35419.5
$ws.Range("L46").Value = 35419.5
$ws.Range("N46").Value = -35731.5
$ws.Range("H70").Value = 10077.667
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10077.667
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10077.667
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -10617.667
$ws.Range("H73").Value = 10077.667
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10077.667
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10077.667
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -11949.667
$ws.Range("H93").Value = 38483.617
$ws.Range("J93").Value = 38483.617
$ws.Range("L93").Value = 38483.617
$ws.Range("N93").Value = -42227.617
$ws.Range("H113").Value = 10362.23
$ws.Range("J113").Value = 2100
$ws.Range("L113").Value = 2100
$ws.Range("N113").Value = -6440
$ws.Range("H122").Value = 2294.889
$ws.Range("I122").Value = 2533.2
$ws.Range("J122").Value = 1997
$ws.Range("K122").Value = 7599.599999999999
$ws.Range("L122").Value = 5991
$ws.Range("M122").Value = -5149.599999999999
$ws.Range("N122").Value = -10891
$ws.Range("H132").Value = 2880.7646
$ws.Range("I132").Value = 2398.2
$ws.Range("K132").Value = 7194.599999999999
$ws.Range("M132").Value = -4664.599999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1571.1111
$ws.Range("I22").Value = 748
$ws.Range("K22").Value = 748
$ws.Range("M22").Value = -453
$ws.Range("H27").Value = 1571.1111
$ws.Range("I27").Value = 748
$ws.Range("K27").Value = 748
$ws.Range("M27").Value = -641
$ws.Range("H55").Value = 527.46155
$ws.Range("I55").Value = 594.73334
$ws.Range("J55").Value = 435.72726
$ws.Range("K55").Value = 594.73334
$ws.Range("L55").Value = 435.72726
$ws.Range("M55").Value = -421.73334
$ws.Range("N55").Value = -781.72726

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6963.0625
$ws.Range("I136").Value = 5913.4346
$ws.Range("K136").Value = 17740.3038
$ws.Range("M136").Value = -15190.3038
